$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the header row (A1:B1)
$ws.Range("A1:B1").HorizontalAlignment = -4108  # xlCenter

# Add the new URL row, matching the style (border) of the existing data rows
$ws.Range("A2:B3").Copy()
$ws.Range("A4:B12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B4").Value = "https://www.facebook.com"
$ws.Range("A4").Value = "URL"

# Resize column B to fit the new, longer URL text
$ws.Columns.Item(2).ColumnWidth = 25.25

$ws.Range("A4").Select()
